$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell (used for the "date" columns,
# which are stored as plain text strings like "2025-10-22" rather than real
# Excel dates) without letting Excel's autoconvert-to-date logic kick in and
# without disturbing the cell's existing style/number-format.
function Set-TextValue($addr, [string]$text) {
    $c = $ws.Range($addr)
    $c.Formula = '="' + $text + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)  # xlPasteValues
}

# --- Row 29 ---
Set-TextValue "N29" "2025-10-22"

# --- Row 30 ---
Set-TextValue "N30" "2025-10-22"
$ws.Range("Q30").Value = 2.29
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.26
$ws.Range("T30").Value = 2.27
$ws.Range("U30").Value = 2.28

# --- Row 42 ---
Set-TextValue "C42" "2025-09-01"
$ws.Range("F42").Value = 4060000
$ws.Range("G42").Value = 4000000
$ws.Range("H42").Value = 4010000
$ws.Range("I42").Value = 3930000
$ws.Range("J42").Value = 4040000

# --- Row 43 ---
Set-TextValue "C43" "2025-09-01"
$ws.Range("F43").Value = 0.04102564102564103

# --- Row 48 ---
Set-TextValue "N48" "2025-10-21"
$ws.Range("Q48").Value = 3.45
$ws.Range("R48").Value = 3.46
$ws.Range("S48").Value = 3.46
$ws.Range("T48").Value = 3.41
$ws.Range("U48").Value = 3.5

# --- Row 49 ---
Set-TextValue "N49" "2025-10-21"
$ws.Range("Q49").Value = 3.56
$ws.Range("R49").Value = 3.58
$ws.Range("S49").Value = 3.59
$ws.Range("T49").Value = 3.55
$ws.Range("U49").Value = 3.63

# --- Row 50 ---
Set-TextValue "N50" "2025-10-21"
$ws.Range("Q50").Value = 3.98
$ws.Range("R50").Value = 4
$ws.Range("S50").Value = 4.02
$ws.Range("T50").Value = 3.99
$ws.Range("U50").Value = 4.05

Write-Output "done"
